# Update algorithm result values in Sheet1 to reflect the newly computed
# RandomForest imputation output ("Update Name of Algo" data refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Cell = "B7";   Value = 4.6984 }
    @{ Cell = "A9";   Value = -21.87080000000002 }
    @{ Cell = "B12";  Value = 5.563599999999996 }
    @{ Cell = "E13";  Value = 16.79680000000001 }
    @{ Cell = "C15";  Value = -13.29289999999999 }
    @{ Cell = "E16";  Value = 16.3332 }
    @{ Cell = "A18";  Value = -22.02270000000001 }
    @{ Cell = "A20";  Value = -19.2933 }
    @{ Cell = "E20";  Value = 16.17079999999999 }
    @{ Cell = "E24";  Value = 16.399 }
    @{ Cell = "B26";  Value = 4.286800000000003 }
    @{ Cell = "A27";  Value = -21.58549999999996 }
    @{ Cell = "B27";  Value = 5.178400000000003 }
    @{ Cell = "B29";  Value = 4.840399999999998 }
    @{ Cell = "B37";  Value = 8.859500000000004 }
    @{ Cell = "B38";  Value = 4.903700000000003 }
    @{ Cell = "C38";  Value = -12.3781 }
    @{ Cell = "E39";  Value = 15.92929999999999 }
    @{ Cell = "C44";  Value = -13.44909999999999 }
    @{ Cell = "E48";  Value = 17.5083 }
    @{ Cell = "B51";  Value = 6.017800000000005 }
    @{ Cell = "C51";  Value = -12.1422 }
    @{ Cell = "E52";  Value = 17.143 }
    @{ Cell = "B55";  Value = 4.919299999999998 }
    @{ Cell = "E56";  Value = 16.536 }
    @{ Cell = "C57";  Value = -13.62779999999999 }
    @{ Cell = "C63";  Value = -11.9095 }
    @{ Cell = "A69";  Value = -21.65199999999999 }
    @{ Cell = "B69";  Value = 5.497399999999995 }
    @{ Cell = "B70";  Value = 5.656400000000004 }
    @{ Cell = "C70";  Value = -11.8828 }
    @{ Cell = "A76";  Value = -19.60609999999999 }
    @{ Cell = "A82";  Value = -21.82450000000001 }
    @{ Cell = "B83";  Value = 5.940999999999998 }
    @{ Cell = "E84";  Value = 16.65669999999999 }
    @{ Cell = "C99";  Value = -12.4213 }
    @{ Cell = "E100"; Value = 16.3843 }
    @{ Cell = "E101"; Value = 16.79570000000001 }
    @{ Cell = "B102"; Value = 8.437600000000005 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
